# Regenerate save_data to use K (strikeouts) instead of Strike# (computed
# from box-score data external to this sheet). Update column G ("K") for
# each game row with the recalculated strikeout totals; std/mean derived
# values downstream of this column are picked up automatically since they
# are not stored separately in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 4
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 3
    18 = 1
    19 = 0
    20 = 2
    21 = 1
    23 = 1
    24 = 3
    25 = 2
    26 = 1
    27 = 0
    28 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 0
    43 = 3
    44 = 1
    45 = 0
    46 = 3
    47 = 2
    48 = 0
    49 = 0
    50 = 0
    51 = 3
    52 = 4
    53 = 0
    54 = 3
    55 = 0
    56 = 3
    57 = 1
    59 = 3
    60 = 0
    61 = 1
    62 = 2
    63 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
